$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '70.940.63'
Set-TextCell $ws.Range("E2") '  +0.12%  '

# Row 3
Set-TextCell $ws.Range("D3") '3.848.97'
Set-TextCell $ws.Range("E3") '  +1.31%  '

# Row 4
Set-TextCell $ws.Range("E4") '  -0.03%  '

# Row 5
Set-TextCell $ws.Range("D5") '707.15'
Set-TextCell $ws.Range("E5") '  +0.80%  '

# Row 6
Set-TextCell $ws.Range("D6") '172.59'
Set-TextCell $ws.Range("E6") '  -0.08%  '

# Row 7
Set-TextCell $ws.Range("D7") '3.847.17'
Set-TextCell $ws.Range("E7") '  +1.29%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.525'
Set-TextCell $ws.Range("E9") '  -0.42%  '

# Row 10
Set-TextCell $ws.Range("E10") '  -0.29%  '

# Row 11
Set-TextCell $ws.Range("D11") '7.33'
Set-TextCell $ws.Range("E11") '  -0.77%  '

# Row 13
Set-TextCell $ws.Range("E13") '  -1.09%  '

# Row 14
Set-TextCell $ws.Range("D14") '36.75'
Set-TextCell $ws.Range("E14") '  +0.88%  '

# Row 15
Set-TextCell $ws.Range("D15") '4.496.35'
Set-TextCell $ws.Range("E15") '  +1.28%  '

# Row 16
Set-TextCell $ws.Range("D16") '3.863.81'
Set-TextCell $ws.Range("E16") '  +1.46%  '

# Row 17
Set-TextCell $ws.Range("D17") '70.957.31'
Set-TextCell $ws.Range("E17") '  +0.18%  '

# Row 18
Set-TextCell $ws.Range("D18") '7.20'
Set-TextCell $ws.Range("E18") '  +0.03%  '

# Row 19
Set-TextCell $ws.Range("E19") '  +0.99%  '

# Row 20
Set-TextCell $ws.Range("D20") '17.36'
Set-TextCell $ws.Range("E20") '  -2.85%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws.Range("D21") '493.38'
Set-TextCell $ws.Range("E21") '  +2.24%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws.Range("D22") '10.66'
Set-TextCell $ws.Range("E22") '  -3.64%  '

# Row 23
Set-TextCell $ws.Range("D23") '0.717'
Set-TextCell $ws.Range("E23") '  +0.42%  '

# Row 24
Set-TextCell $ws.Range("D24") '85.31'
Set-TextCell $ws.Range("E24") '  +1.16%  '

# Row 25
Set-TextCell $ws.Range("E25") '  +1.95%  '

# Row 26
Set-TextCell $ws.Range("D26") '10.67'
Set-TextCell $ws.Range("E26") '  +1.89%  '

# Row 27
Set-TextCell $ws.Range("D27") '12.16'
Set-TextCell $ws.Range("E27") '  -2.23%  '

# Row 28
Set-TextCell $ws.Range("E28") '  -3.29%  '

# Row 29
Set-TextCell $ws.Range("E29") '  +1.90%  '

# Row 30
Set-TextCell $ws.Range("E30") '  -0.06%  '

# Row 31
Set-TextCell $ws.Range("D31") '7.50'
Set-TextCell $ws.Range("E31") '  -0.37%  '

# Row 32
Set-TextCell $ws.Range("E32") '  -0.95%  '

# Row 33
Set-TextCell $ws.Range("D33") '29.44'
Set-TextCell $ws.Range("E33") '  -0.22%  '

# Row 34
Set-TextCell $ws.Range("E34") '  -1.69%  '

# Row 35
Set-TextCell $ws.Range("D35") '3.803.83'
Set-TextCell $ws.Range("E35") '  +1.48%  '

# Row 36
Set-TextCell $ws.Range("D36") '9.16'
Set-TextCell $ws.Range("E36") '  -0.72%  '

# Row 37
Set-TextCell $ws.Range("E37") '  -0.01%  '

# Row 38
Set-TextCell $ws.Range("E38") '  +0.21%  '

# Row 39
Set-TextCell $ws.Range("E39") '  +6.93%  '

# Row 40
Set-TextCell $ws.Range("E40") '  +6.41%  '

# Row 41
Set-TextCell $ws.Range("D41") '6.04'
Set-TextCell $ws.Range("E41") '  -0.24%  '

# Row 42
Set-TextCell $ws.Range("E42") '  -3.48%  '

# Row 44
Set-TextCell $ws.Range("E44") '  +0.15%  '

# Row 45
Set-TextCell $ws.Range("D45") '0.000318'
Set-TextCell $ws.Range("E45") '  -3.06%  '

# Row 46
Set-TextCell $ws.Range("D46") '163.13'
Set-TextCell $ws.Range("E46") '  +0.65%  '

# Row 47
Set-TextCell $ws.Range("D47") '48.70'
Set-TextCell $ws.Range("E47") '  +0.02%  '

# Row 48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell $ws.Range("D48") '1.39'
Set-TextCell $ws.Range("E48") '  +0.48%  '

# Row 49
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws.Range("D49") '416.28'
Set-TextCell $ws.Range("E49") '  +1.74%  '

# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws.Range("D50") '8.63'
Set-TextCell $ws.Range("E50") '  +0.64%  '

# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws.Range("D51") '0.299'
Set-TextCell $ws.Range("E51") '  -1.44%  '
